$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: clone a cell's value + hyperlink (if any) from a source cell to a
# destination cell on the same worksheet.
# ---------------------------------------------------------------------------
function Copy-CellWithHyperlink($ws, $srcAddr, $dstAddr) {
    $srcRange = $ws.Range($srcAddr)
    $dstRange = $ws.Range($dstAddr)
    $srcAddress = $srcRange.Address()

    $linkUrl = $null
    $linkDisplay = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $srcAddress) {
            $linkUrl = $h.Address
            $linkDisplay = $h.TextToDisplay
        }
    }

    $dstRange.Value = $srcRange.Value()

    if ($linkUrl -ne $null) {
        $ws.Hyperlinks.Add($dstRange, $linkUrl, "", "", $linkDisplay) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This shared string is used on the Overview sheet (Status columns for
#    zh-cn/de-de) as well as on the per-language sheets (File Extension
#    column, which carries the same text in the source data).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$handedBackText = "Handed back: in sync with en-US"

$wsOverview.Range("B2").Value = $handedBackText
$wsOverview.Range("C2").Value = $handedBackText
$wsOverview.Range("B3").Value = $handedBackText
$wsOverview.Range("C3").Value = $handedBackText

$wsZh.Range("C2").Value = $handedBackText
$wsZh.Range("C3").Value = $handedBackText

$wsDe.Range("C2").Value = $handedBackText
$wsDe.Range("C3").Value = $handedBackText

# ---------------------------------------------------------------------------
# 2. Generate the handback report columns: F (Latest Target File) and
#    G (Latest Handback File) mirror the existing A (Source File Name) and
#    D (Latest Handoff File) columns respectively, for rows 2 and 3, on
#    both the zh-cn and de-de sheets.
# ---------------------------------------------------------------------------
Copy-CellWithHyperlink $wsZh "A2" "F2"
Copy-CellWithHyperlink $wsZh "D2" "G2"
Copy-CellWithHyperlink $wsZh "A3" "F3"
Copy-CellWithHyperlink $wsZh "D3" "G3"

Copy-CellWithHyperlink $wsDe "A2" "F2"
Copy-CellWithHyperlink $wsDe "D2" "G2"
Copy-CellWithHyperlink $wsDe "A3" "F3"
Copy-CellWithHyperlink $wsDe "D3" "G3"

# ---------------------------------------------------------------------------
# 3. Latest Handback DateTime (column H): record the handback timestamps.
#    zh-cn rows share one timestamp, de-de rows share a later one.
# ---------------------------------------------------------------------------
$wsZh.Range("H2").Value = "2016-03-20 18:48:59"
$wsZh.Range("H3").Value = "2016-03-20 18:48:59"

$wsDe.Range("H2").Value = "2016-03-20 18:49:06"
$wsDe.Range("H3").Value = "2016-03-20 18:49:06"
